# Weekly update: a new price record for "Jengibre" (Mercado Mayorista Lo
# Valledor de Santiago) is inserted at row 44, pushing the existing
# historical rows (old rows 44-83) down by one (new rows 45-84).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 44 - this shifts rows 44:83 down to
# 45:84 (and grows the used range from R83 to R84), matching the diff.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new weekly record.
$ws.Cells.Item(44, 1).Value  = 6
$ws.Cells.Item(44, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(44, 3).Value  = "Metropolitana"
$ws.Cells.Item(44, 4).Value  = 44740
$ws.Cells.Item(44, 5).Value  = 13
$ws.Cells.Item(44, 6).Value  = 100114007
$ws.Cells.Item(44, 7).Value  = "Jengibre"
$ws.Cells.Item(44, 8).Value  = "Sin especificar"
$ws.Cells.Item(44, 9).Value  = "Primera"
$ws.Cells.Item(44, 10).Value = 400
$ws.Cells.Item(44, 11).Value = 13000
$ws.Cells.Item(44, 12).Value = 14000
$ws.Cells.Item(44, 13).Value = 13425
$ws.Cells.Item(44, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(44, 15).Value = "Perú"
$ws.Cells.Item(44, 16).Value = 1033
$ws.Cells.Item(44, 17).Value = 13
$ws.Cells.Item(44, 18).Value = "Hortaliza"
